# Auto-generated: applies scheduled market-data refresh values
# to the Pandaemonium_Profits workbook (per-cell updates across all job sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 313.54
$ws.Range("I15").Value = 313.54
$ws.Range("K15").Value = 940.6200000000001
$ws.Range("M15").Value = -771.6200000000001

$ws.Range("H132").Value = 1365.6428
$ws.Range("I132").Value = 1455.1333
$ws.Range("J132").Value = 999.5454999999999
$ws.Range("K132").Value = 4365.3999
$ws.Range("L132").Value = 2998.6365
$ws.Range("M132").Value = -1835.3999
$ws.Range("N132").Value = -8058.6365

$ws.Range("H135").Value = 31915194
$ws.Range("I135").Value = 11905023
$ws.Range("J135").Value = 200000620
$ws.Range("K135").Value = 107145207
$ws.Range("L135").Value = 1800005580
$ws.Range("M135").Value = -107142672
$ws.Range("N135").Value = -1800010650

$ws.Range("H137").Value = 1735.3939
$ws.Range("I137").Value = 1185.8695
$ws.Range("J137").Value = 2999.3
$ws.Range("K137").Value = 3557.6085
$ws.Range("L137").Value = 8997.900000000001
$ws.Range("M137").Value = -1007.6085
$ws.Range("N137").Value = -14097.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6208
$ws.Range("I61").Value = 3169.2188
$ws.Range("K61").Value = 3169.2188
$ws.Range("M61").Value = -2957.2188

$ws.Range("H74").Value = 6027.6665
$ws.Range("I74").Value = 1997.1578
$ws.Range("J74").Value = 21343.6
$ws.Range("K74").Value = 1997.1578
$ws.Range("L74").Value = 21343.6
$ws.Range("M74").Value = -1123.1578
$ws.Range("N74").Value = -23091.6

$ws.Range("H77").Value = 6027.6665
$ws.Range("I77").Value = 1997.1578
$ws.Range("J77").Value = 21343.6
$ws.Range("K77").Value = 9985.789000000001
$ws.Range("L77").Value = 106718
$ws.Range("M77").Value = -5617.789000000001
$ws.Range("N77").Value = -115454

$ws.Range("H132").Value = 9745.833000000001
$ws.Range("I132").Value = 4615.1665
$ws.Range("J132").Value = 12311.167
$ws.Range("K132").Value = 13845.4995
$ws.Range("L132").Value = 36933.501
$ws.Range("M132").Value = -11315.4995
$ws.Range("N132").Value = -41993.501

$ws.Range("H136").Value = 6208
$ws.Range("I136").Value = 3169.2188
$ws.Range("K136").Value = 9507.6564
$ws.Range("M136").Value = -6957.6564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -826

$ws.Range("H134").Value = 35223.582
$ws.Range("I134").Value = 2807.9524
$ws.Range("J134").Value = 103296.4
$ws.Range("K134").Value = 8423.8572
$ws.Range("L134").Value = 309889.2
$ws.Range("M134").Value = -5888.8572
$ws.Range("N134").Value = -314959.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5168.647
$ws.Range("I31").Value = 5074.4443
$ws.Range("J31").Value = 5532
$ws.Range("K31").Value = 5074.4443
$ws.Range("L31").Value = 5532
$ws.Range("M31").Value = -4779.4443
$ws.Range("N31").Value = -6122

$ws.Range("H34").Value = 5168.647
$ws.Range("I34").Value = 5074.4443
$ws.Range("J34").Value = 5532
$ws.Range("K34").Value = 5074.4443
$ws.Range("L34").Value = 5532
$ws.Range("M34").Value = -4872.4443
$ws.Range("N34").Value = -5936

$ws.Range("H58").Value = 1247053.8
$ws.Range("I58").Value = 1716363.5
$ws.Range("J58").Value = 3382.9
$ws.Range("K58").Value = 1716363.5
$ws.Range("L58").Value = 3382.9
$ws.Range("M58").Value = -1716160.5
$ws.Range("N58").Value = -3788.9

$ws.Range("H134").Value = 2153.1184
$ws.Range("I134").Value = 1235
$ws.Range("J134").Value = 4026.08
$ws.Range("K134").Value = 3705
$ws.Range("L134").Value = 12078.24
$ws.Range("M134").Value = -1170
$ws.Range("N134").Value = -17148.24

$ws.Range("H136").Value = 1247053.8
$ws.Range("I136").Value = 1716363.5
$ws.Range("J136").Value = 3382.9
$ws.Range("K136").Value = 5149090.5
$ws.Range("L136").Value = 10148.7
$ws.Range("M136").Value = -5146540.5
$ws.Range("N136").Value = -15248.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5054129
$ws.Range("I5").Value = 406.79166
$ws.Range("J5").Value = 18530722
$ws.Range("K5").Value = 1220.37498
$ws.Range("L5").Value = 55592166
$ws.Range("M5").Value = -1108.37498
$ws.Range("N5").Value = -55592390

$ws.Range("H68").Value = 6397.1577
$ws.Range("I68").Value = 810
$ws.Range("J68").Value = 8392.571
$ws.Range("K68").Value = 2430
$ws.Range("L68").Value = 25177.713
$ws.Range("M68").Value = -1619
$ws.Range("N68").Value = -26799.713

$ws.Range("H71").Value = 6397.1577
$ws.Range("I71").Value = 810
$ws.Range("J71").Value = 8392.571
$ws.Range("K71").Value = 7290
$ws.Range("L71").Value = 75533.139
$ws.Range("M71").Value = -3234
$ws.Range("N71").Value = -83645.139

$ws.Range("H107").Value = 1712.8235
$ws.Range("I107").Value = 389.75
$ws.Range("J107").Value = 2888.889
$ws.Range("K107").Value = 1169.25
$ws.Range("L107").Value = 8666.667000000001
$ws.Range("M107").Value = 750.75
$ws.Range("N107").Value = -12506.667

$ws.Range("H122").Value = 700.2143
$ws.Range("I122").Value = 446.05884
$ws.Range("J122").Value = 873.04
$ws.Range("K122").Value = 4014.52956
$ws.Range("L122").Value = 7857.36
$ws.Range("M122").Value = -1564.52956
$ws.Range("N122").Value = -12757.36

$ws.Range("H135").Value = 5054129
$ws.Range("I135").Value = 406.79166
$ws.Range("J135").Value = 18530722
$ws.Range("K135").Value = 3661.12494
$ws.Range("L135").Value = 166776498
$ws.Range("M135").Value = -1126.12494
$ws.Range("N135").Value = -166781568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5737.7144
$ws.Range("I70").Value = 5224
$ws.Range("J70").Value = 6368.1816
$ws.Range("K70").Value = 5224
$ws.Range("L70").Value = 6368.1816
$ws.Range("M70").Value = -4954
$ws.Range("N70").Value = -6908.1816

$ws.Range("H73").Value = 5737.7144
$ws.Range("I73").Value = 5224
$ws.Range("J73").Value = 6368.1816
$ws.Range("K73").Value = 5224
$ws.Range("L73").Value = 6368.1816
$ws.Range("M73").Value = -4288
$ws.Range("N73").Value = -8240.1816

$ws.Range("H102").Value = 7354.3335
$ws.Range("I102").Value = 7937.3335
$ws.Range("J102").Value = 6771.3335
$ws.Range("K102").Value = 7937.3335
$ws.Range("L102").Value = 6771.3335
$ws.Range("M102").Value = -6315.3335
$ws.Range("N102").Value = -10015.3335

$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 1880
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 5640
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -3170
$ws.Range("N126").Value = -16190

$ws.Range("H132").Value = 8999.129000000001
$ws.Range("I132").Value = 7337.476
$ws.Range("J132").Value = 12488.6
$ws.Range("K132").Value = 22012.428
$ws.Range("L132").Value = 37465.8
$ws.Range("M132").Value = -19482.428
$ws.Range("N132").Value = -42525.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 273
$ws.Range("I55").Value = 268.22223
$ws.Range("K55").Value = 268.22223
$ws.Range("M55").Value = -95.22223000000002

$ws.Range("H68").Value = 1367.9231
$ws.Range("I68").Value = 974.75
$ws.Range("K68").Value = 974.75
$ws.Range("M68").Value = -225.75

$ws.Range("H71").Value = 1367.9231
$ws.Range("I71").Value = 974.75
$ws.Range("K71").Value = 4873.75
$ws.Range("M71").Value = -1129.75

$ws.Range("H122").Value = 8633.333000000001
$ws.Range("I122").Value = 8200
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 24600
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -22150
$ws.Range("N122").Value = -33400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3699.5293
$ws.Range("I122").Value = 2849.0715
$ws.Range("J122").Value = 7668.3335
$ws.Range("K122").Value = 8547.2145
$ws.Range("L122").Value = 23005.0005
$ws.Range("M122").Value = -6097.2145
$ws.Range("N122").Value = -27905.0005

$ws.Range("H126").Value = 2099.375
$ws.Range("I126").Value = 2071.4285
$ws.Range("J126").Value = 2295
$ws.Range("K126").Value = 6214.2855
$ws.Range("L126").Value = 6885
$ws.Range("M126").Value = -3744.2855
$ws.Range("N126").Value = -11825

$ws.Range("H132").Value = 1115.8143
$ws.Range("I132").Value = 443.44446
$ws.Range("J132").Value = 2326.08
$ws.Range("K132").Value = 1330.33338
$ws.Range("L132").Value = 6978.24
$ws.Range("M132").Value = 1199.66662
$ws.Range("N132").Value = -12038.24

$ws.Range("H136").Value = 3765.137
$ws.Range("I136").Value = 2563.373
$ws.Range("K136").Value = 7690.119000000001
$ws.Range("M136").Value = -5140.119000000001
